$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Angpt1"
$ws.Range("C2").Value = "Itga5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 20.22494433333334
$ws.Range("H2").Value = 60.67483300000001
$ws.Range("I2").Value = 0.9541201174409912
$ws.Range("J2").Value = 0.9683117878380343
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 33.211442
$ws.Range("N2").Value = 99.634326
$ws.Range("O2").Value = 0.211580186305583
$ws.Range("P2").Value = 0.2175281749633597
$ws.Range("Q2").Value = 671.6995656797287
$ws.Range("R2").Value = 6045.296091117559
$ws.Range("S2").Value = 0.2018729122060696
$ws.Range("T2").Value = 0.2106350960039156

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Angpt1"
$ws.Range("C3").Value = "Itga5"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 20.22494433333334
$ws.Range("H3").Value = 60.67483300000001
$ws.Range("I3").Value = 0.9541201174409912
$ws.Range("J3").Value = 0.9683117878380343
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 51.17424933333334
$ws.Range("N3").Value = 153.522748
$ws.Range("O3").Value = 0.3260158715178649
$ws.Range("P3").Value = 0.3351809012869699
$ws.Range("Q3").Value = 1034.996344066787
$ws.Range("R3").Value = 9314.967096601085
$ws.Range("S3").Value = 0.3110583016202524
$ws.Range("T3").Value = 0.3245596177743496

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Angpt1"
$ws.Range("C4").Value = "Itga5"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 20.22494433333334
$ws.Range("H4").Value = 60.67483300000001
$ws.Range("I4").Value = 0.9541201174409912
$ws.Range("J4").Value = 0.9683117878380343
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.19450366666666
$ws.Range("N4").Value = 66.58351099999999
$ws.Range("O4").Value = 0.1413945597650736
$ws.Range("P4").Value = 0.1453694746776606
$ws.Range("Q4").Value = 448.8826011642959
$ws.Range("R4").Value = 4039.943410478663
$ws.Range("S4").Value = 0.1349073939685693
$ws.Range("T4").Value = 0.1407629759222014

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Angpt1"
$ws.Range("C5").Value = "Itga5"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 20.22494433333334
$ws.Range("H5").Value = 60.67483300000001
$ws.Range("I5").Value = 0.9541201174409912
$ws.Range("J5").Value = 0.9683117878380343
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 37.51216133333333
$ws.Range("N5").Value = 112.536484
$ws.Range("O5").Value = 0.2389787857941174
$ws.Range("P5").Value = 0.2456970098971044
$ws.Range("Q5").Value = 758.6813747896858
$ws.Range("R5").Value = 6828.132373107172
$ws.Range("S5").Value = 0.2280144671677887
$ws.Range("T5").Value = 0.2379113109199244

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Angpt1"
$ws.Range("C6").Value = "Itga5"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 20.22494433333334
$ws.Range("H6").Value = 60.67483300000001
$ws.Range("I6").Value = 0.9541201174409912
$ws.Range("J6").Value = 0.9683117878380343
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 12.8762265
$ws.Range("N6").Value = 25.752453
$ws.Range("O6").Value = 0.08203059661736112
$ws.Range("P6").Value = 0.05622443917490542
$ws.Range("Q6").Value = 260.4209641858916
$ws.Range("R6").Value = 1562.525785115349
$ws.Range("S6").Value = 0.07826704247831116
$ws.Range("T6").Value = 0.05444278721764348

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Angpt1"
$ws.Range("C7").Value = "Itga5"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.040521
$ws.Range("H7").Value = 0.121563
$ws.Range("I7").Value = 0.001911594941455862
$ws.Range("J7").Value = 0.00194002818046413
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 33.211442
$ws.Range("N7").Value = 99.634326
$ws.Range("O7").Value = 0.211580186305583
$ws.Range("P7").Value = 0.2175281749633597
$ws.Range("Q7").Value = 1.345760841282
$ws.Range("R7").Value = 12.111847571538
$ws.Range("S7").Value = 0.0004044556138540413
$ws.Range("T7").Value = 0.0004220107894738496

# Row 8
$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "Angpt1"
$ws.Range("C8").Value = "Itga5"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.040521
$ws.Range("H8").Value = 0.121563
$ws.Range("I8").Value = 0.001911594941455862
$ws.Range("J8").Value = 0.00194002818046413
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 51.17424933333334
$ws.Range("N8").Value = 153.522748
$ws.Range("O8").Value = 0.3260158715178649
$ws.Range("P8").Value = 0.3351809012869699
$ws.Range("Q8").Value = 2.073631757236
$ws.Range("R8").Value = 18.662685815124
$ws.Range("S8").Value = 0.0006232102908278748
$ws.Range("T8").Value = 0.0006502603940500874

# Row 9
$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "Angpt1"
$ws.Range("C9").Value = "Itga5"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.040521
$ws.Range("H9").Value = 0.121563
$ws.Range("I9").Value = 0.001911594941455862
$ws.Range("J9").Value = 0.00194002818046413
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 22.19450366666666
$ws.Range("N9").Value = 66.58351099999999
$ws.Range("O9").Value = 0.1413945597650736
$ws.Range("P9").Value = 0.1453694746776606
$ws.Range("Q9").Value = 0.8993434830769999
$ws.Range("R9").Value = 8.094091347692999
$ws.Range("S9").Value = 0.0002702891251962933
$ws.Range("T9").Value = 0.0002820208774539283

# Row 10
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Angpt1"
$ws.Range("C10").Value = "Itga5"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.040521
$ws.Range("H10").Value = 0.121563
$ws.Range("I10").Value = 0.001911594941455862
$ws.Range("J10").Value = 0.00194002818046413
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 37.51216133333333
$ws.Range("N10").Value = 112.536484
$ws.Range("O10").Value = 0.2389787857941174
$ws.Range("P10").Value = 0.2456970098971044
$ws.Range("Q10").Value = 1.520030289388
$ws.Range("R10").Value = 13.680272604492
$ws.Range("S10").Value = 0.0004568306380392988
$ws.Range("T10").Value = 0.0004766591230561568

# Row 11
$ws.Range("A11").Value = "M1"
$ws.Range("B11").Value = "Angpt1"
$ws.Range("C11").Value = "Itga5"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.040521
$ws.Range("H11").Value = 0.121563
$ws.Range("I11").Value = 0.001911594941455862
$ws.Range("J11").Value = 0.00194002818046413
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 12.8762265
$ws.Range("N11").Value = 25.752453
$ws.Range("O11").Value = 0.08203059661736112
$ws.Range("P11").Value = 0.05622443917490542
$ws.Range("Q11").Value = 0.5217575740065
$ws.Range("R11").Value = 3.130545444039
$ws.Range("S11").Value = 0.0001568092735383539
$ws.Range("T11").Value = 0.0001090769964301079

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Angpt1"
$ws.Range("C12").Value = "Itga5"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.9320170000000001
$ws.Range("H12").Value = 1.864034
$ws.Range("I12").Value = 0.04396828761755308
$ws.Range("J12").Value = 0.02974818398150157
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 33.211442
$ws.Range("N12").Value = 99.634326
$ws.Range("O12").Value = 0.211580186305583
$ws.Range("P12").Value = 0.2175281749633597
$ws.Range("Q12").Value = 30.953628538514
$ws.Range("R12").Value = 185.721771231084
$ws.Range("S12").Value = 0.009302818485659338
$ws.Range("T12").Value = 0.006471068169970286

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Angpt1"
$ws.Range("C13").Value = "Itga5"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.9320170000000001
$ws.Range("H13").Value = 1.864034
$ws.Range("I13").Value = 0.04396828761755308
$ws.Range("J13").Value = 0.02974818398150157
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 51.17424933333334
$ws.Range("N13").Value = 153.522748
$ws.Range("O13").Value = 0.3260158715178649
$ws.Range("P13").Value = 0.3351809012869699
$ws.Range("Q13").Value = 47.69527034090534
$ws.Range("R13").Value = 286.1716220454321
$ws.Range("S13").Value = 0.01433435960678472
$ws.Range("T13").Value = 0.009971023118570297

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Angpt1"
$ws.Range("C14").Value = "Itga5"
$ws.Range("D14").Value = "M1"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.9320170000000001
$ws.Range("H14").Value = 1.864034
$ws.Range("I14").Value = 0.04396828761755308
$ws.Range("J14").Value = 0.02974818398150157
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 22.19450366666666
$ws.Range("N14").Value = 66.58351099999999
$ws.Range("O14").Value = 0.1413945597650736
$ws.Range("P14").Value = 0.1453694746776606
$ws.Range("Q14").Value = 20.68565472389567
$ws.Range("R14").Value = 124.113928343374
$ws.Range("S14").Value = 0.006216876671308057
$ws.Range("T14").Value = 0.004324477878005281

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Angpt1"
$ws.Range("C15").Value = "Itga5"
$ws.Range("D15").Value = "M2"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.9320170000000001
$ws.Range("H15").Value = 1.864034
$ws.Range("I15").Value = 0.04396828761755308
$ws.Range("J15").Value = 0.02974818398150157
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 37.51216133333333
$ws.Range("N15").Value = 112.536484
$ws.Range("O15").Value = 0.2389787857941174
$ws.Range("P15").Value = 0.2456970098971044
$ws.Range("Q15").Value = 34.96197206940933
$ws.Range("R15").Value = 209.771832416456
$ws.Range("S15").Value = 0.01050748798828936
$ws.Range("T15").Value = 0.007309039854123873

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Angpt1"
$ws.Range("C16").Value = "Itga5"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.9320170000000001
$ws.Range("H16").Value = 1.864034
$ws.Range("I16").Value = 0.04396828761755308
$ws.Range("J16").Value = 0.02974818398150157
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 12.8762265
$ws.Range("N16").Value = 25.752453
$ws.Range("O16").Value = 0.08203059661736112
$ws.Range("P16").Value = 0.05622443917490542
$ws.Range("Q16").Value = 12.0008619938505
$ws.Range("R16").Value = 48.00344797540201
$ws.Range("S16").Value = 0.003606744865511611
$ws.Range("T16").Value = 0.00167257496083183
